# Kilimanjaro Weekly Scoreboard - append this week's workout rows (178-184)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array matches the sheet's column order:
# A Participant, B Date, C Workout Type, D Total Duration, E Total Distance,
# F Total Elevation, G Zone 1, H Zone 2, I Zone 3, J Zone 4, K Zone 5,
# L Workout Level, M Week
$newRows = @(
    @("Matt",   45481, "Walk",    50, 1.95, 56, 50,  0,  0, 0, 0, "Sauntering Hippo", 5),
    @("Steven", 45481, "Workout", 27,    0,  0, 27,  0,  0, 0, 0, "Brave Leopard",    5),
    @("Matt",   45481, "Walk",    45, 1.73, 59, 45,  0,  0, 0, 0, "Sauntering Hippo", 5),
    @("Steven", 45481, "Workout", 35,    0,  0,  1, 10, 21, 3, 0, "Brave Leopard",    5),
    @("Eric",   45481, "Workout", 66,    0,  0, 23, 40,  4, 1, 0, "Wily Hyena",       5),
    @("Steven", 45481, "Walk",    25, 1.26, 33, 25,  0,  0, 0, 0, "Brave Leopard",    5),
    @("Steven", 45482, "Walk",    26, 1.28, 23, 26,  0,  0, 0, 0, "Brave Leopard",    5)
)

$startRow = 178
$endRow = $startRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $values[$c]
    }
}

# Reuse the existing date number format (style already used by column B above)
# instead of letting Excel mint a brand-new style for the new date cells.
$ws.Range("B177").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Land the final selection where Excel would leave it after typing the
# last new row (the sheet's used range/dimension grows automatically to
# A1:M184 as a side effect of the cell writes above).
$ws.Range("N" + $endRow).Select()
